$wb = $excel.ActiveWorkbook

function Set-F($SheetIndex, $Row, $Value) {
    $ws = $wb.Worksheets.Item($SheetIndex)
    $ws.Cells.Item($Row, 6).Value = $Value
}

# Sheet 1: 展览
Set-F 1 4  7827
Set-F 1 5  7827
Set-F 1 8  2123
Set-F 1 9  8606
Set-F 1 13 5749
Set-F 1 15 2727
Set-F 1 18 414
Set-F 1 22 43
Set-F 1 23 3853
Set-F 1 25 58
Set-F 1 26 50
Set-F 1 28 17
Set-F 1 29 5356
Set-F 1 35 380
Set-F 1 36 1982
Set-F 1 40 4063
Set-F 1 44 3485
Set-F 1 50 8

# Sheet 2: 演出
Set-F 2 2 120
Set-F 2 6 17
Set-F 2 8 32

# Sheet 3: 本地生活
Set-F 3 2 269
Set-F 3 3 1360

# Sheet 4: 全部类型
Set-F 4 2  269
Set-F 4 3  1360
Set-F 4 4  7827
Set-F 4 5  7827
Set-F 4 8  2123
Set-F 4 9  8606
Set-F 4 13 5749
Set-F 4 15 2727
Set-F 4 18 414
Set-F 4 20 120
Set-F 4 25 43
Set-F 4 26 3853
Set-F 4 28 58
Set-F 4 29 50
Set-F 4 30 17
Set-F 4 31 5356
Set-F 4 35 380
Set-F 4 36 1982
Set-F 4 38 17
Set-F 4 42 4063
Set-F 4 46 3485
Set-F 4 47 32
